# Apply updates to "BB Alert Result Statistics" workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update data values in B4:D11 ---
$data = @{
    4  = @(21,  33.35,   8.8309523809523807)
    5  = @(91,  586.27,  55.85230769230769)
    6  = @(103, 629.01,  64.342233009708764)
    7  = @(338, 507.5,   43.209171597633173)
    8  = @(272, 1028.21, 43.639301470588229)
    9  = @(198, 318.35000000000002, 38.082373737373729)
    10 = @(107, 518.17999999999995, 36.638504672897191)
    11 = @(95,  215.38,  39.774947368421053)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
}

# --- Set column widths for B, C, D (bestFit / customWidth) ---
# (values chosen so the engine's internal pixel-rounding lands on the
#  width closest to the target stored width)
$ws.Range("B:B").ColumnWidth = 4.75
$ws.Range("C:C").ColumnWidth = 7.25
$ws.Range("D:D").ColumnWidth = 11.25

# --- Update zoom level ---
$win = $excel.ActiveWindow
$win.Zoom = 220

# --- Update selection ---
$ws.Range("A6:D6").Select()
